# Add a fourth "SUBSCRIPTION" block (rows 16-19) to Sheet1, mirroring the
# existing ROLE/ROUTE/... access-control blocks (rows 1-4, 6-9, 11-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Merge the destination B column FIRST (on still-empty/default-styled
#        cells) so that when we paste formats afterwards they cleanly
#        overwrite whatever the merge operation touched, instead of the
#        merge clobbering formats we already set. ---
$ws.Range("B17:B19").Merge()

# --- 2. Clone formatting (styles/borders/fills/merge look) from the third
#        block (rows 11-14: header + STUDENT/TEACHER/ADMIN) onto the new
#        block (rows 16-19). Formats-only paste so shared styles are reused
#        instead of Excel cloning new style records. ---
$src = $ws.Range("A11:J14")
$dst = $ws.Range("A16:J19")
$src.Copy()
$dst.PasteSpecial(-4122)

# --- 3. Fill in the header row (row 16) - identical header to every other
#        block. ---
$ws.Range("A16").Value2 = "ROLE"
$ws.Range("B16").Value2 = "ROUTE"
$ws.Range("C16").Value2 = "CREATE 1"
$ws.Range("D16").Value2 = "DELETE 1"
$ws.Range("E16").Value2 = "UPDATE 1"
$ws.Range("F16").Value2 = "READ 1"
$ws.Range("G16").Value2 = "BULK CREATE"
$ws.Range("H16").Value2 = "BULK DELETE"
$ws.Range("I16").Value2 = "BULK UPDATE"
$ws.Range("J16").Value2 = "BULK READ"

# --- 4. Fill in the data rows (17-19): STUDENT / TEACHER / ADMIN, with the
#        merged B column labelled with the new "SUBSCRIPTION" row/model
#        name. ---
$ws.Range("A17").Value2 = "STUDENT"
$ws.Range("B17").Value2 = "SUBSCRIPTION"
$ws.Range("A18").Value2 = "TEACHER"
$ws.Range("A19").Value2 = "ADMIN"

# --- 5. Widen column B slightly so the new "SUBSCRIPTION" label fits. ---
$ws.Columns.Item(2).ColumnWidth = 13.83

# --- 6. Match the saved cursor/selection position recorded in the sheet. ---
$ws.Range("H24").Select() | Out-Null
